$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new value, whether the value is numeric-looking text
# (numeric-looking text must be forced to Text format so Excel keeps it as a string
# and preserves exact formatting such as trailing zeros).
$changes = @(
    @('D2', '244.97', 1),
    @('D3', '21.96', 1),
    @('D4', '5.400', 1),
    @('D5', '0.06006', 1),
    @('D6', '3.386', 1),
    @('D7', '0.8120', 1),
    @('D8', '0.9575', 1),
    @('B9', 'One', 0),
    @('C9', 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one', 0),
    @('D9', '0.01118', 1),
    @('E9', '8OneONEBestin24h', 0),
    @('B10', 'WazirX', 0),
    @('C10', 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx', 0),
    @('D10', '0.1424', 1),
    @('E10', '9WazirXWRX', 0),
    @('B11', 'MandalaExchangeToken', 0),
    @('C11', 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx', 0),
    @('D11', '0.07387', 1),
    @('E11', '10MandalaExchangeTokenMDX', 0),
    @('B12', 'LiechtensteinCryptoassetsExchange', 0),
    @('C12', 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx', 0),
    @('D12', '0.03389', 1),
    @('E12', '11LiechtensteinCryptoassetsExchangeLCX', 0),
    @('B13', 'BitrueCoin', 0),
    @('C13', 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr', 0),
    @('D13', '0.03057', 1),
    @('E13', '12BitrueCoinBTR', 0),
    @('B14', 'BitMartToken', 0),
    @('C14', 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx', 0),
    @('D14', '0.09419', 1),
    @('E14', '13BitMartTokenBMX', 0),
    @('B15', 'MCDex', 0),
    @('C15', 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb', 0),
    @('D15', '4.002', 1),
    @('E15', '14MCDexMCB', 0),
    @('B16', 'BitForexToken', 0),
    @('C16', 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf', 0),
    @('D16', '0.001589', 1),
    @('E16', '15BitForexTokenBF', 0),
    @('B17', 'CoinExToken', 0),
    @('C17', 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet', 0),
    @('D17', '0.04802', 1),
    @('E17', '16CoinExTokenCET', 0),
    @('D18', '0.006215', 1),
    @('D20', '0.0009886', 1),
    @('D22', '3.697', 1),
    @('D23', '6.401', 1),
    @('D26', '0.1284', 1),
    @('D40', '0.04031', 1),
    @('B41', 'KickToken', 0),
    @('C41', 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick', 0),
    @('D41', '0.006493', 1),
    @('E41', '40KickTokenKICK', 0),
    @('B42', 'BKEXToken', 0),
    @('C42', 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk', 0),
    @('D42', '0.1071', 1),
    @('E42', '41BKEXTokenBKK', 0),
    @('B43', 'CEJI', 0),
    @('C43', 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji', 0),
    @('D43', '0.002901', 1),
    @('E43', '42CEJICEJI', 0),
    @('D44', '0.005843', 1),
    @('D45', '0.00005259', 1),
    @('E47', '46CoinbaseStockTokenCOIN', 0),
    @('D48', '0.02207', 1)
)

foreach ($item in $changes) {
    $cellRef = $item[0]
    $newVal = $item[1]
    $isNumericText = $item[2]
    $rng = $ws.Range($cellRef)
    if ($isNumericText -eq 1) {
        $rng.NumberFormat = "@"
    }
    $rng.Value2 = $newVal
}

Write-Output "Applied $($changes.Count) cell updates"
